# C5-PowerPoint.pptx edit:
#  1. Slide 6's table is re-styled from the "No Style, Table Grid" table
#     style to a different built-in table style.
#  2. The deck's design theme's colour scheme is switched from the
#     "Integral" palette to the standard Office "Office Theme" palette
#     (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{40669362-D528-49BE-8AE1-833DB88B590B}")
    }
}

# --- 2. Swap the design theme's colour scheme -----------------------------
# Integral  -> Office Theme
#   dk1      000000 -> 000000
#   lt1      FFFFFF -> FFFFFF
#   dk2      455F51 -> 44546A
#   lt2      E3DED1 -> E7E6E6
#   accent1  99CB38 -> 5B9BD5
#   accent2  63A537 -> ED7D31
#   accent3  E6D024 -> A5A5A5
#   accent4  CC9700 -> FFC000
#   accent5  4EB3CF -> 4472C4
#   accent6  378DA6 -> 70AD47
#   hlink    6B9F25 -> 0563C1
#   folHlink B26B02 -> 954F72
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0          # dk1      #000000
$colors.Colors(2).RGB  = 16777215   # lt1      #FFFFFF
$colors.Colors(3).RGB  = 6968388    # dk2      #44546A
$colors.Colors(4).RGB  = 15132391   # lt2      #E7E6E6
$colors.Colors(5).RGB  = 13998939   # accent1  #5B9BD5
$colors.Colors(6).RGB  = 3243501    # accent2  #ED7D31
$colors.Colors(7).RGB  = 10855845   # accent3  #A5A5A5
$colors.Colors(8).RGB  = 49407      # accent4  #FFC000
$colors.Colors(9).RGB  = 12874308   # accent5  #4472C4
$colors.Colors(10).RGB = 4697456    # accent6  #70AD47
$colors.Colors(11).RGB = 12673797   # hlink    #0563C1
$colors.Colors(12).RGB = 7491477    # folHlink #954F72
